# Update a handful of statistic values on the Brazil Summary sheet.
# Each cell currently holds its value as text (not a number), so we
# explicitly force a text number-format before writing the new value;
# otherwise Excel would auto-convert the numeric-looking string into a
# numeric cell, which would change the cell's underlying type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): Micro 22.9 -> 22.95
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "22.95"

# Enterprises density (per 1000 people): MSMEs 25.9 -> 25.94
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "25.94"

# Employment (% of total): MSMEs 53 -> 53.01
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.01"

# Enterprises (% of total): SMEs 11.5 -> 11.51
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "11.51"

# Enterprises (% of total): MSMEs 99.6 -> 99.61
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "99.61"
